$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.038.04"
$ws.Range("E2").Value = "  -0.90%  "

$ws.Range("D3").Value = "1.784.49"
$ws.Range("E3").Value = "  -2.26%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "227.03"
$ws.Range("E5").Value = "  -1.76%  "

$ws.Range("D6").Value = "0.552"

$ws.Range("D8").Value = "31.26"
$ws.Range("E8").Value = "  -1.05%  "

$ws.Range("D9").Value = "46.19"
$ws.Range("E9").Value = "  +0.84%  "

$ws.Range("E10").Value = "  -0.72%  "

$ws.Range("D11").Value = "0.0658"
$ws.Range("E11").Value = "  -2.85%  "

$ws.Range("E12").Value = "  -0.30%  "

$ws.Range("D13").Value = "2.044.91"
$ws.Range("E13").Value = "  -2.17%  "

$ws.Range("D14").Value = "11.39"
$ws.Range("E14").Value = "  +10.72%  "

$ws.Range("D15").Value = "1.786.05"
$ws.Range("E15").Value = "  -2.33%  "

$ws.Range("D16").Value = "0.633"
$ws.Range("E16").Value = "  -2.20%  "

$ws.Range("D17").Value = "34.049.03"
$ws.Range("E17").Value = "  -0.91%  "

$ws.Range("D18").Value = "4.22"
$ws.Range("E18").Value = "  -3.14%  "

$ws.Range("D19").Value = "69.37"
$ws.Range("E19").Value = "  -0.94%  "

$ws.Range("D20").Value = "252.66"
$ws.Range("E20").Value = "  -2.81%  "

$ws.Range("D21").Value = "0.0₃0741"
$ws.Range("E21").Value = "  -1.65%  "

$ws.Range("E22").Value = "  +0.17%  "

$ws.Range("D23").Value = "10.43"
$ws.Range("E23").Value = "  -1.16%  "

$ws.Range("D24").Value = "4.24"
$ws.Range("E24").Value = "  -3.55%  "

$ws.Range("E25").Value = "  -2.63%  "

$ws.Range("D26").Value = "157.03"
$ws.Range("E26").Value = "  -2.71%  "

$ws.Range("D27").Value = "16.55"
$ws.Range("E27").Value = "  -1.67%  "

$ws.Range("D28").Value = "7.02"
$ws.Range("E28").Value = "  -2.54%  "

$ws.Range("E29").Value = "  -2.17%  "

$ws.Range("E30").Value = "  +0.10%  "

$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.80"
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = "  -2.54%  "

$ws.Range("D32").Value = "0.0516"
$ws.Range("E32").Value = "  -0.47%  "

$ws.Range("E33").Value = "  -1.19%  "

$ws.Range("D34").Value = "3.59"
$ws.Range("E34").Value = "  +0.45%  "

$ws.Range("D35").Value = "1.84"
$ws.Range("E35").Value = "  -0.64%  "

$ws.Range("D36").Value = "1.449.29"
$ws.Range("E36").Value = "  -8.14%  "

$ws.Range("E37").Value = "  -1.30%  "

$ws.Range("D38").Value = "0.626"
$ws.Range("E38").Value = "  -1.00%  "

$ws.Range("D39").Value = "0.0186"
$ws.Range("E39").Value = "  -1.72%  "

$ws.Range("D40").Value = "83.38"
$ws.Range("E40").Value = "  -2.04%  "

$ws.Range("D41").Value = "2.83"
$ws.Range("E41").Value = "  -1.31%  "

$ws.Range("E42").Value = "  -0.52%  "

$ws.Range("D43").Value = "0.898"
$ws.Range("E43").Value = "  -2.43%  "

$ws.Range("D44").Value = "2.07"
$ws.Range("E44").Value = "  -3.23%  "

$ws.Range("D45").Value = "0.0509"
$ws.Range("E45").Value = "  -2.26%  "

$ws.Range("E46").Value = "  +0.03%  "

$ws.Range("D47").Value = "1.943.38"
$ws.Range("E47").Value = "  -1.98%  "

$ws.Range("D48").Value = "5.76"
$ws.Range("E48").Value = "  +0.05%  "

$ws.Range("E49").Value = "  +0.12%  "

$ws.Range("D50").Value = "11.88"
$ws.Range("E50").Value = "  +5.30%  "

$ws.Range("D51").Value = "51.13"
$ws.Range("E51").Value = "  -3.43%  "
